$wb = $excel.ActiveWorkbook

# --- Sheets: rename Sheet1, add Sheet2 and Sheet3 in order ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Customer Registration Data"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Sheet3"

# --- Populate "Customer Registration Data" sheet ---
$ws1.Range("A1").Value = "first name"
$ws1.Range("B1").Value = "last name"
$ws1.Range("C1").Value = "email"
$ws1.Range("D1").Value = "telephone"
$ws1.Range("E1").Value = "password"

$ws1.Range("A2").Value = "mike 1 "
$ws1.Range("B2").Value = "pence"
$ws1.Range("C2").Value = "mike1@gmail.com"
$ws1.Range("D2").Value = 7474847575
$ws1.Range("E2").Value = "teyf6464646"

$ws1.Range("A3").Value = "mike 2"
$ws1.Range("B3").Value = "pence"
$ws1.Range("C3").Value = "mike2@gmail.com"
$ws1.Range("D3").Value = 7474847575
$ws1.Range("E3").Value = "teyf6464646"

$ws1.Range("A4").Value = "mike 3"
$ws1.Range("B4").Value = "pence"
$ws1.Range("C4").Value = "mike3@gmail.com"
$ws1.Range("D4").Value = 7474847575
$ws1.Range("E4").Value = "teyf6464646"

$ws1.Range("A5").Value = "mike 4"
$ws1.Range("B5").Value = "pence"
$ws1.Range("C5").Value = "mike4@gmail.com"
$ws1.Range("D5").Value = 7474847575
$ws1.Range("E5").Value = "teyf6464646"

# --- Header row formatting: built-in "Heading 1" style + bottom border ---
$headerRange = $ws1.Range("A1:E1")
$headerRange.Style = "Heading 1"

# --- Telephone column number format ---
$ws1.Range("D1:D5").NumberFormat = "[<=9999999]###\-####;\(###\)\ ###\-####"

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 17.5703125
$ws1.Columns.Item(2).ColumnWidth = 16.28515625
$ws1.Columns.Item(3).ColumnWidth = 20.5703125
$ws1.Columns.Item(4).ColumnWidth = 29.28515625
$ws1.Columns.Item(5).ColumnWidth = 17.85546875

# --- Sheet view tweaks ---
$ws1.Range("D16").Select()
$excel.ActiveWindow.Zoom = 120

$ws3.Range("K15").Select()

$wb.Windows.Item(1).WindowState = -4143
